$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(1.02, 1.05548994593055, 1.061424838910845, 1.060929646222222, 1.069241641411597, 1.027664988806262, 1.060496535907449, 1.064148935661988, 1.063655088705775, 1.071944673888349, 1.023883598764913)
  3 = @(1.02, 1.057481369175367, 1.063275134532518, 1.062738332709983, 1.071255413547409, 1.027887981298487, 1.062134628253618, 1.065811242061347, 1.065275791089603, 1.073771610744368, 1.024427228538382)
  4 = @(1.02, 1.058758038831842, 1.064459948300748, 1.063896466363182, 1.072541087766212, 1.028022304815294, 1.06318251664347, 1.066874009047761, 1.066311872155894, 1.074936025470244, 1.024774928114871)
  5 = @(1.02, 1.059291944954545, 1.064955107618312, 1.06438046666693, 1.073077480529279, 1.028076398096406, 1.063620201019465, 1.067317760689389, 1.066744460236651, 1.075421350637182, 1.024920141427052)
  6 = @(1.02, 1.059381426926763, 1.065038076112889, 1.06446156484451, 1.073167304057407, 1.028085341517596, 1.063693524307367, 1.067392091727687, 1.066816920073142, 1.075502594325947, 1.024944467496168)
  7 = @(1.02, 1.058765183875222, 1.064466576123321, 1.063902944854679, 1.072548271130285, 1.028023036936786, 1.063188376140011, 1.066879950351214, 1.06631766407569, 1.074942526824398, 1.024776872214848)
  8 = @(1.02, 1.056165452285753, 1.062052760461774, 1.061543454649135, 1.069925835211832, 1.027742420294643, 1.061052660627605, 1.064713407932098, 1.064205451515515, 1.072565799603614, 1.024068171122109)
  9 = @(1.02, 1.051490962055854, 1.057701855708943, 1.057290194516348, 1.065169176649201, 1.02717110012795, 1.057194911940072, 1.06079522733845, 1.060384850022335, 1.068239447463576, 1.022787572975079)
  10 = @(1.02, 1.048308631709967, 1.054732665886663, 1.054387460938941, 1.061903155163125, 1.026737782867706, 1.054556822031751, 1.058112643404826, 1.057768624512969, 1.06525868794841, 1.021911536755081)
  11 = @(1.02, 1.046914243763213, 1.053429983208215, 1.053113893776385, 1.060465529960267, 1.026537518923452, 1.053398117789737, 1.056933651240332, 1.05661869156694, 1.063944243894346, 1.021526691607142)
  12 = @(1.02, 1.046393772277612, 1.052943489233188, 1.052638267099776, 1.059927933704339, 1.026461215333608, 1.052965200467847, 1.056493042371218, 1.056188926895142, 1.063452355230138, 1.021382894398438)
  13 = @(1.02, 1.04650553094353, 1.053047963461358, 1.052740407855249, 1.060043414083565, 1.026477669754926, 1.053058177865684, 1.056587676715177, 1.056281232822316, 1.063558033260848, 1.021413778117695)
  14 = @(1.02, 1.046871273513582, 1.053389823290499, 1.053074631073138, 1.060421165966752, 1.02653125085515, 1.053362384533514, 1.056897285388112, 1.056583221064873, 1.06390365904992, 1.021514822714324)
  15 = @(1.02, 1.047096281783841, 1.053600105295297, 1.053280214911967, 1.060653432104353, 1.026564009392258, 1.053549480143407, 1.05708768865025, 1.056768935837562, 1.064116125072206, 1.021576966568677)
  16 = @(1.02, 1.048400821075373, 1.054818756911619, 1.054471626942195, 1.061998065496052, 1.026750805984388, 1.05463337095659, 1.058190517060474, 1.057844576581393, 1.065345416312222, 1.021936959770857)
  17 = @(1.02, 1.049214682077722, 1.055578587469832, 1.055214462761095, 1.062835193817646, 1.026864583282488, 1.055308836405776, 1.058877586505081, 1.058514680055693, 1.066110105247074, 1.022161283503912)
  18 = @(1.02, 1.049687814024209, 1.056020147691605, 1.05564614289431, 1.063321222882128, 1.026929729906088, 1.055701246392921, 1.059276666482247, 1.058903895386263, 1.066553848671739, 1.022291596950322)
  19 = @(1.02, 1.049848873886686, 1.056170432732814, 1.055793064575306, 1.063486566336034, 1.026951737207495, 1.055834782249306, 1.059412460003218, 1.059036330588156, 1.066704768115439, 1.022335941026589)
  20 = @(1.02, 1.049127526336431, 1.055497234531145, 1.055134929808812, 1.062745611552766, 1.026852502141947, 1.05523652895772, 1.058804044254156, 1.058442954843661, 1.066028298427547, 1.022137270733933)
  21 = @(1.02, 1.046763641938444, 1.053289226941841, 1.052976282054661, 1.06031002744381, 1.026515525619684, 1.053272873393122, 1.056806187808599, 1.0564943661119, 1.063801982139431, 1.021485091188423)
  22 = @(1.02, 1.045262679392997, 1.051885774325423, 1.051604168451994, 1.058757818298359, 1.026292555932877, 1.052023616415002, 1.055534525300513, 1.055253971979251, 1.062381075742562, 1.021070119871744)
  23 = @(1.02, 1.046059784066538, 1.052631233451916, 1.052332984795865, 1.059582678842084, 1.026411814944775, 1.052687278536646, 1.05621015129679, 1.055912993877556, 1.063136354799227, 1.021290577270125)
  24 = @(1.02, 1.049166913141186, 1.055533999474427, 1.055170872277755, 1.062786096881938, 1.026857964851273, 1.055269206452606, 1.058837280007313, 1.05847536945846, 1.066065270469714, 1.022148122710162)
  25 = @(1.02, 1.052710820976279, 1.058838524072846, 1.05840138514778, 1.066415329671975, 1.027327981266478, 1.058203695290075, 1.061820367550788, 1.061384538619968, 1.069374653135731, 1.023122497920252)
}

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item($row, $cols[$i]).Value = $vals[$i]
  }
}
